# Generate Report for Handback
# Rename the "zh-tw" locale sheet/table to "ru-ru" and refresh the
# handoff/handback timestamps captured in the report.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the locale worksheet ------------------------------------
$wsLocale = $wb.Worksheets.Item("zh-tw")
$wsLocale.Name = "ru-ru"

# --- 2. Update the Overview sheet's column header -----------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B1").Value = "ru-ru"

# --- 3. Rename the locale table to match -------------------------------
$localeTable = $wsLocale.ListObjects.Item(1)
$localeTable.Name = "ru-ru"

# --- 4. Refresh the "Correspond Handoff Datetime" column (E2:E5) --------
$wsLocale.Range("E2").Value = "2016-03-11 01:04:17"
$wsLocale.Range("E3").Value = "2016-03-11 01:04:17"
$wsLocale.Range("E4").Value = "2016-03-11 01:04:17"
$wsLocale.Range("E5").Value = "2016-03-11 01:04:17"

# --- 5. Refresh the "Correspond Handback DateTime" column (H2:H5) -------
$wsLocale.Range("H2").Value = "2016-03-16 23:59:51"
$wsLocale.Range("H3").Value = "2016-03-16 23:59:51"
$wsLocale.Range("H4").Value = "2016-03-17 16:45:12"
$wsLocale.Range("H5").Value = "2016-03-16 23:59:51"
